$wb = $excel.ActiveWorkbook
$wsP = $wb.Worksheets.Item("Processes")
$wb.Names.Add("TestName", "=Processes!`$D`$2:`$D`$8")
$wsP.Columns("D:D").Cut()
$wsP.Columns("B:B").Insert(-4161)
Write-Host $wb.Names.Item("TestName").RefersTo
